# Recompute the "IPC PO" (predicted output) column using a zero-weight model:
# with weight = 0, the predicted output is always 0, so:
#   IPC PO (C)   = 0
#   DELTA  (D)   = IPC PO - IPC RO = 0 - B = -B
#   DELTA^2 (E)  = (IPC RO)^2  (since DELTA = -B, DELTA^2 = B^2)
# Row 52 (TOTAL) sums the DELTA/DELTA^2 columns, row 53 (MSE) is the mean of DELTA^2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 51
$totalRow = 52
$mseRow = 53

$colRO = 2   # B: IPC RO
$colPO = 3   # C: IPC PO
$colDelta = 4   # D: DELTA
$colDelta2 = 5   # E: DELTA^2

$sumDelta = 0
$sumDelta2 = 0
$count = 0

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ipcRO = $ws.Cells.Item($r, $colRO).Value2

    # Weight reset to 0 => predicted output is 0 for every point.
    $ipcPO = 0
    $delta = $ipcPO - $ipcRO
    $delta2 = $delta * $delta

    $ws.Cells.Item($r, $colPO).Value2 = $ipcPO
    $ws.Cells.Item($r, $colDelta).Value2 = $delta
    $ws.Cells.Item($r, $colDelta2).Value2 = $delta2

    $sumDelta = $sumDelta + $delta
    $sumDelta2 = $sumDelta2 + $delta2
    $count = $count + 1
}

$ws.Cells.Item($totalRow, $colPO).Value2 = $sumDelta
$ws.Cells.Item($totalRow, $colDelta2).Value2 = $sumDelta2

$ws.Cells.Item($mseRow, $colDelta2).Value2 = $sumDelta2 / $count
